$d = $word.ActiveDocument

$d.Content.Find.Execute("45×43=", $true, $false, $false, $false, $false, $true, 1, $false, "88×20=", 2) | Out-Null
$d.Content.Find.Execute("52×76=", $true, $false, $false, $false, $false, $true, 1, $false, "47×63=", 2) | Out-Null
$d.Content.Find.Execute("93×43=", $true, $false, $false, $false, $false, $true, 1, $false, "66×60=", 2) | Out-Null
$d.Content.Find.Execute("72×19=", $true, $false, $false, $false, $false, $true, 1, $false, "19×57=", 2) | Out-Null
$d.Content.Find.Execute("71×68=", $true, $false, $false, $false, $false, $true, 1, $false, "60×85=", 2) | Out-Null
$d.Content.Find.Execute("65×96=", $true, $false, $false, $false, $false, $true, 1, $false, "80×26=", 2) | Out-Null
$d.Content.Find.Execute("39×12=", $true, $false, $false, $false, $false, $true, 1, $false, "21×11=", 2) | Out-Null
$d.Content.Find.Execute("43×97=", $true, $false, $false, $false, $false, $true, 1, $false, "17×35=", 2) | Out-Null
$d.Content.Find.Execute("73×23=", $true, $false, $false, $false, $false, $true, 1, $false, "91×60=", 2) | Out-Null
$d.Content.Find.Execute("11×70=", $true, $false, $false, $false, $false, $true, 1, $false, "61×22=", 2) | Out-Null
$d.Content.Find.Execute("70×43=", $true, $false, $false, $false, $false, $true, 1, $false, "95×85=", 2) | Out-Null
$d.Content.Find.Execute("45×20=", $true, $false, $false, $false, $false, $true, 1, $false, "34×66=", 2) | Out-Null
$d.Content.Find.Execute("45×18=", $true, $false, $false, $false, $false, $true, 1, $false, "23×88=", 2) | Out-Null
$d.Content.Find.Execute("90×62=", $true, $false, $false, $false, $false, $true, 1, $false, "91×85=", 2) | Out-Null
$d.Content.Find.Execute("17×87=", $true, $false, $false, $false, $false, $true, 1, $false, "47×12=", 2) | Out-Null
$d.Content.Find.Execute("89×46=", $true, $false, $false, $false, $false, $true, 1, $false, "75×33=", 2) | Out-Null
$d.Content.Find.Execute("26×66=", $true, $false, $false, $false, $false, $true, 1, $false, "21×41=", 2) | Out-Null
$d.Content.Find.Execute("22×22=", $true, $false, $false, $false, $false, $true, 1, $false, "64×26=", 2) | Out-Null
$d.Content.Find.Execute("37×50=", $true, $false, $false, $false, $false, $true, 1, $false, "45×71=", 2) | Out-Null
$d.Content.Find.Execute("90×64=", $true, $false, $false, $false, $false, $true, 1, $false, "99×42=", 2) | Out-Null
$d.Content.Find.Execute("17×11=", $true, $false, $false, $false, $false, $true, 1, $false, "86×85=", 2) | Out-Null
$d.Content.Find.Execute("95×24=", $true, $false, $false, $false, $false, $true, 1, $false, "37×45=", 2) | Out-Null
$d.Content.Find.Execute("32×27=", $true, $false, $false, $false, $false, $true, 1, $false, "58×19=", 2) | Out-Null
$d.Content.Find.Execute("19×42=", $true, $false, $false, $false, $false, $true, 1, $false, "92×13=", 2) | Out-Null
$d.Content.Find.Execute("26×72=", $true, $false, $false, $false, $false, $true, 1, $false, "16×46=", 2) | Out-Null
